# Update LR-pairs data (Efna3-Epha7) with new TPM-derived values.
# The sending/target cluster grid changes from {ECs,FAPs} x {ECs,FAPs,MuSCs,Resolving-Mac}
# to {ECs,FAPs,MuSCs} x {ECs,FAPs,MuSCs} (Resolving-Mac cluster replaced by MuSCs as a
# sending cluster too), and all numeric measurements (columns E:T) are refreshed.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Sending cluster / Ligand symbol / Receptor symbol / Target cluster for rows 2..10
$labels = @(
    @("ECs",  "Efna3", "Epha7", "ECs"),
    @("ECs",  "Efna3", "Epha7", "FAPs"),
    @("ECs",  "Efna3", "Epha7", "MuSCs"),
    @("FAPs", "Efna3", "Epha7", "ECs"),
    @("FAPs", "Efna3", "Epha7", "FAPs"),
    @("FAPs", "Efna3", "Epha7", "MuSCs"),
    @("MuSCs","Efna3", "Epha7", "ECs"),
    @("MuSCs","Efna3", "Epha7", "FAPs"),
    @("MuSCs","Efna3", "Epha7", "MuSCs")
)

# Numeric data columns E..T for rows 2..10
$numbers = @(
    @(2, 0.6666666666666666, 0.2420556666666667, 0.726167, 0.5314769098578004, 0.5314769098578004, 1, 0.3333333333333333, 0.2282403333333334, 0.684721, 0.1514399067192403, 0.1514399067192403, 0.05524686604522223, 0.497221794407, 0.08048681365229539, 0.08048681365229537),
    @(2, 0.6666666666666666, 0.2420556666666667, 0.726167, 0.5314769098578004, 0.5314769098578004, 3, 1, 1.151356333333333, 3.454069, 0.763937263734893, 0.763937263734893, 0.2786923248358889, 2.508230923523, 0.4060150162550444, 0.4060150162550444),
    @(2, 0.6666666666666666, 0.2420556666666667, 0.726167, 0.5314769098578004, 0.5314769098578004, 3, 1, 0.127538, 0.3826140000000001, 0.08462282954586674, 0.08462282954586674, 0.03087129561533334, 0.277841660538, 0.04497507995046062, 0.04497507995046062),
    @(1, 0.3333333333333333, 0.2054156666666667, 0.616247, 0.4510271759376837, 0.4510271759376837, 1, 0.3333333333333333, 0.2282403333333334, 0.684721, 0.1514399067192403, 0.1514399067192403, 0.04688414023188889, 0.421957262087, 0.06830351345184521, 0.0683035134518452),
    @(1, 0.3333333333333333, 0.2054156666666667, 0.616247, 0.4510271759376837, 0.4510271759376837, 3, 1, 1.151356333333333, 3.454069, 0.763937263734893, 0.763937263734893, 0.2365066287825556, 2.128559659043, 0.3445564666559102, 0.3445564666559102),
    @(1, 0.3333333333333333, 0.2054156666666667, 0.616247, 0.4510271759376837, 0.4510271759376837, 3, 1, 0.127538, 0.3826140000000001, 0.08462282954586674, 0.08462282954586674, 0.02619830329533334, 0.235784729658, 0.03816719582992825, 0.03816719582992825),
    @(1, 0.3333333333333333, 0.007968333333333332, 0.023905, 0.01749591420451593, 0.01749591420451593, 1, 0.3333333333333333, 0.2282403333333334, 0.684721, 0.1514399067192403, 0.1514399067192403, 0.001818695056111111, 0.016368255505, 0.002649579615099724, 0.002649579615099723),
    @(1, 0.3333333333333333, 0.007968333333333332, 0.023905, 0.01749591420451593, 0.01749591420451593, 3, 1, 1.151356333333333, 3.454069, 0.763937263734893, 0.763937263734893, 0.009174391049444444, 0.082569519445, 0.01336578082393834, 0.01336578082393834),
    @(1, 0.3333333333333333, 0.007968333333333332, 0.023905, 0.01749591420451593, 0.01749591420451593, 3, 1, 0.127538, 0.3826140000000001, 0.08462282954586674, 0.08462282954586674, 0.001016265296666667, 0.009146387670000001, 0.00148055376547786, 0.00148055376547786)
)

$cols = @("E","F","G","H","I","J","K","L","M","N","O","P","Q","R","S","T")

for ($i = 0; $i -lt $labels.Length; $i++) {
    $row = 2 + $i
    $lbl = $labels[$i]
    $ws.Range("A$row").Value = $lbl[0]
    $ws.Range("B$row").Value = $lbl[1]
    $ws.Range("C$row").Value = $lbl[2]
    $ws.Range("D$row").Value = $lbl[3]

    $nums = $numbers[$i]
    for ($j = 0; $j -lt $cols.Length; $j++) {
        $ws.Range("$($cols[$j])$row").Value = $nums[$j]
    }
}
